$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: rename unitary_weight/total_weight columns to
#     unitary_measure/total_measure (new shared strings, old ones become
#     orphaned and will be dropped from sharedStrings.xml on save).
$ws.Range("C1").Value = "unitary_measure"
$ws.Range("D1").Value = "total_measure"

# --- Per-row cleanup: each data row kept either its "unitary" (C) or its
#     "total" (D) measure, not both (the redundant one is cleared).
$ws.Range("C2").ClearContents()
$ws.Range("D3").ClearContents()
$ws.Range("D4").ClearContents()
$ws.Range("D5").ClearContents()
$ws.Range("C6").ClearContents()
$ws.Range("D7").ClearContents()
$ws.Range("C9").ClearContents()
$ws.Range("C10").ClearContents()
$ws.Range("C11").ClearContents()
$ws.Range("D12").ClearContents()
$ws.Range("C13").ClearContents()
$ws.Range("C14").ClearContents()
$ws.Range("D16").ClearContents()
$ws.Range("D17").ClearContents()
$ws.Range("D18").ClearContents()

# --- Column widths: column A gets its own best-fit width (~53.66 chars),
#     B:E keep the original uniform width untouched.
$ws.Columns.Item(1).ColumnWidth = 52.75

# --- View state: move the selection from A17 to C2 (this also clears the
#     stored topLeftCell scroll anchor, resetting the scrolled-down view).
$ws.Range("C2").Select()
